$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '55.137.12'
$ws.Range("E2").Value = '  -4.31%  '

$ws.Range("D3").Value = '2.926.70'
$ws.Range("E3").Value = '  -7.18%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '''477.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -9.28%  '

$ws.Range("D6").Value = '''129.42'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.21%  '

$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("D8").Value = '2.925.50'
$ws.Range("E8").Value = '  -7.08%  '

$ws.Range("D9").Value = '''0.411'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.54%  '

$ws.Range("D10").Value = '''6.75'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.26%  '

$ws.Range("D11").Value = '''0.0988'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -11.07%  '

$ws.Range("D12").Value = '''0.340'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -12.71%  '

$ws.Range("E13").Value = '  -2.02%  '

$ws.Range("D14").Value = '3.422.66'
$ws.Range("E14").Value = '  -7.13%  '

$ws.Range("D15").Value = '''23.78'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -7.99%  '

$ws.Range("D16").Value = '54.999.90'
$ws.Range("E16").Value = '  -4.53%  '

$ws.Range("D17").Value = '2.915.87'
$ws.Range("E17").Value = '  -7.27%  '

$ws.Range("E18").Value = '  -10.91%  '

$ws.Range("D19").Value = '''5.49'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.47%  '

$ws.Range("D20").Value = '''11.71'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -11.06%  '

$ws.Range("D21").Value = '''7.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -10.22%  '

$ws.Range("D22").Value = '''306.88'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -11.85%  '

$ws.Range("E23").Value = '  +0.00%  '

$ws.Range("D24").Value = '''0.451'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -11.85%  '

$ws.Range("D25").Value = '''59.53'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -14.25%  '

$ws.Range("D26").Value = '''0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.07%  '

$ws.Range("E27").Value = '  -7.06%  '

$ws.Range("E28").Value = '  +0.01%  '

$ws.Range("D29").Value = '0.0₃0824'
$ws.Range("E29").Value = '  -14.61%  '

$ws.Range("D30").Value = '''6.35'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.63%  '

$ws.Range("D31").Value = '''6.37'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.15%  '

$ws.Range("D32").Value = '''1.15'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.63%  '

$ws.Range("E33").Value = '  -12.13%  '

$ws.Range("D34").Value = '''19.05'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -12.22%  '

$ws.Range("D35").Value = '''145.51'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.73%  '

$ws.Range("D36").Value = '''4.27'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -13.73%  '

$ws.Range("D37").Value = '''5.49'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -12.27%  '

$ws.Range("D38").Value = '''1.24'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -12.32%  '

$ws.Range("D39").Value = '''23.38'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -9.92%  '

$ws.Range("D40").Value = '''0.0630'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -9.57%  '

$ws.Range("D41").Value = '2.949.22'
$ws.Range("E41").Value = '  -7.05%  '

$ws.Range("D42").Value = '''0.999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.01%  '

$ws.Range("D43").Value = '''35.74'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -11.54%  '

$ws.Range("B44").Value = 'ONDO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D44").Value = '''0.978'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -9.43%  '

$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = '''0.619'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -10.82%  '

$ws.Range("D46").Value = '''1.35'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.82%  '

$ws.Range("D47").Value = '''3.47'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -12.42%  '

$ws.Range("D48").Value = '2.073.02'
$ws.Range("E48").Value = '  -8.54%  '

$ws.Range("D49").Value = '''5.45'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -12.54%  '

$ws.Range("D50").Value = '''0.0221'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.05%  '

$ws.Range("D51").Value = '''18.27'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -10.87%  '
